# Auto-generated edit script: applies per-cell numeric updates
# described by the Moogle_Profits.xlsx diff (scheduled-runner price refresh).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 835263.5600000001
$ws.Range("I15").Value = 835263.5600000001
$ws.Range("K15").Value = 2505790.68
$ws.Range("M15").Value = -2505621.68
$ws.Range("H17").Value = 2040.375
$ws.Range("J17").Value = 2040.375
$ws.Range("L17").Value = 6121.125
$ws.Range("N17").Value = -6457.125
$ws.Range("H86").Value = 6665.3794
$ws.Range("I86").Value = 6048.222
$ws.Range("J86").Value = 6943.1
$ws.Range("K86").Value = 6048.222
$ws.Range("L86").Value = 6943.1
$ws.Range("M86").Value = -4925.222
$ws.Range("N86").Value = -9189.1
$ws.Range("H89").Value = 6665.3794
$ws.Range("I89").Value = 6048.222
$ws.Range("J89").Value = 6943.1
$ws.Range("K89").Value = 30241.11
$ws.Range("L89").Value = 34715.5
$ws.Range("M89").Value = -24625.11
$ws.Range("N89").Value = -45947.5
$ws.Range("H92").Value = 812.3570999999999
$ws.Range("I92").Value = 874.36365
$ws.Range("J92").Value = 585
$ws.Range("K92").Value = 874.36365
$ws.Range("L92").Value = 585
$ws.Range("M92").Value = 373.63635
$ws.Range("N92").Value = -3081
$ws.Range("H106").Value = 44005028
$ws.Range("I106").Value = 44005028
$ws.Range("K106").Value = 44005028
$ws.Range("M106").Value = -44004397
$ws.Range("H116").Value = 9987.111000000001
$ws.Range("I116").Value = 9978.4
$ws.Range("J116").Value = 9998
$ws.Range("K116").Value = 9978.4
$ws.Range("L116").Value = 9998
$ws.Range("M116").Value = -6536.4
$ws.Range("N116").Value = -16882
$ws.Range("H132").Value = 2205.0625
$ws.Range("I132").Value = 2122.2666
$ws.Range("K132").Value = 6366.7998
$ws.Range("M132").Value = -3836.7998
$ws.Range("H138").Value = 3157.6191
$ws.Range("I138").Value = 2334
$ws.Range("K138").Value = 7002
$ws.Range("M138").Value = -1862

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1954.3334
$ws.Range("I45").Value = 1545.2
$ws.Range("J45").Value = 4000
$ws.Range("K45").Value = 1545.2
$ws.Range("L45").Value = 4000
$ws.Range("M45").Value = -1168.2
$ws.Range("N45").Value = -4754
$ws.Range("H74").Value = 2543.7058
$ws.Range("I74").Value = 738
$ws.Range("K74").Value = 738
$ws.Range("M74").Value = 136
$ws.Range("H77").Value = 2543.7058
$ws.Range("I77").Value = 738
$ws.Range("K77").Value = 3690
$ws.Range("M77").Value = 678
$ws.Range("H122").Value = 2641.6428
$ws.Range("J122").Value = 4329.6665
$ws.Range("L122").Value = 12988.9995
$ws.Range("N122").Value = -17888.9995

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6762.364
$ws.Range("I20").Value = 6986.375
$ws.Range("J20").Value = 6165
$ws.Range("K20").Value = 6986.375
$ws.Range("L20").Value = 6165
$ws.Range("M20").Value = -6739.375
$ws.Range("N20").Value = -6659
$ws.Range("H80").Value = 16107.077
$ws.Range("I80").Value = 50514.5
$ws.Range("J80").Value = 814.8889
$ws.Range("K80").Value = 50514.5
$ws.Range("L80").Value = 814.8889
$ws.Range("M80").Value = -49516.5
$ws.Range("N80").Value = -2810.8889
$ws.Range("H83").Value = 16107.077
$ws.Range("I83").Value = 50514.5
$ws.Range("J83").Value = 814.8889
$ws.Range("K83").Value = 252572.5
$ws.Range("L83").Value = 4074.4445
$ws.Range("M83").Value = -247580.5
$ws.Range("N83").Value = -14058.4445
$ws.Range("H107").Value = 4224.75
$ws.Range("I107").Value = 4224.75
$ws.Range("K107").Value = 4224.75
$ws.Range("M107").Value = -2304.75
$ws.Range("H134").Value = 3014.7646
$ws.Range("I134").Value = 1518.3572
$ws.Range("K134").Value = 4555.071599999999
$ws.Range("M134").Value = -2020.071599999999

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 3125.76
$ws.Range("I105").Value = 2807.75
$ws.Range("J105").Value = 3691.111
$ws.Range("K105").Value = 2807.75
$ws.Range("L105").Value = 3691.111
$ws.Range("M105").Value = -1060.75
$ws.Range("N105").Value = -7185.111
$ws.Range("H122").Value = 2233.087
$ws.Range("I122").Value = 2165.2856
$ws.Range("K122").Value = 6495.8568
$ws.Range("M122").Value = -4045.8568
$ws.Range("H132").Value = 3952.348
$ws.Range("I132").Value = 3115.2632
$ws.Range("K132").Value = 9345.7896
$ws.Range("M132").Value = -6815.7896
$ws.Range("H134").Value = 4695.793
$ws.Range("I134").Value = 2842.5
$ws.Range("J134").Value = 13591.6
$ws.Range("K134").Value = 8527.5
$ws.Range("L134").Value = 40774.8
$ws.Range("M134").Value = -5992.5
$ws.Range("N134").Value = -45844.8

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 189.71428
$ws.Range("I5").Value = 220.25
$ws.Range("J5").Value = 149
$ws.Range("K5").Value = 660.75
$ws.Range("L5").Value = 447
$ws.Range("M5").Value = -548.75
$ws.Range("N5").Value = -671
$ws.Range("H135").Value = 189.71428
$ws.Range("I135").Value = 220.25
$ws.Range("J135").Value = 149
$ws.Range("K135").Value = 1982.25
$ws.Range("L135").Value = 1341
$ws.Range("M135").Value = 552.75
$ws.Range("N135").Value = -6411

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 14864.75
$ws.Range("H97").Value = 759.2857
$ws.Range("J97").Value = 766.2
$ws.Range("L97").Value = 766.2
$ws.Range("N97").Value = -1758.2
$ws.Range("H99").Value = 104826.55
$ws.Range("I99").Value = 3485.5
$ws.Range("J99").Value = 162735.72
$ws.Range("K99").Value = 3485.5
$ws.Range("L99").Value = 162735.72
$ws.Range("M99").Value = -1239.5
$ws.Range("N99").Value = -167227.72
$ws.Range("H107").Value = 330.375
$ws.Range("J107").Value = 521.75
$ws.Range("L107").Value = 521.75
$ws.Range("N107").Value = -4361.75
$ws.Range("H122").Value = 3016.9656
$ws.Range("I122").Value = 1250.2273
$ws.Range("J122").Value = 8569.571
$ws.Range("K122").Value = 3750.6819
$ws.Range("L122").Value = 25708.713
$ws.Range("M122").Value = -1300.6819
$ws.Range("N122").Value = -30608.713
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").Value = ""
$ws.Range("H132").Value = 3048.311
$ws.Range("I132").Value = 2222.4102
$ws.Range("J132").Value = 8416.666999999999
$ws.Range("K132").Value = 6667.230599999999
$ws.Range("L132").Value = 25250.001
$ws.Range("M132").Value = -4137.230599999999
$ws.Range("N132").Value = -30310.001

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 300000
$ws.Range("J43").Value = 300000
$ws.Range("L43").Value = 300000
$ws.Range("N43").Value = -300386
$ws.Range("H68").Value = 5097.278
$ws.Range("I68").Value = 4620.125
$ws.Range("K68").Value = 4620.125
$ws.Range("M68").Value = -3871.125
$ws.Range("H71").Value = 5097.278
$ws.Range("I71").Value = 4620.125
$ws.Range("K71").Value = 23100.625
$ws.Range("M71").Value = -19356.625
$ws.Range("H81").Value = 181805.56
$ws.Range("J81").Value = 181805.56
$ws.Range("L81").Value = 181805.56
$ws.Range("N81").Value = -183801.56
$ws.Range("H84").Value = 181805.56
$ws.Range("J84").Value = 181805.56
$ws.Range("L84").Value = 545416.6799999999
$ws.Range("N84").Value = -555400.6799999999
$ws.Range("H132").Value = 11227.481
$ws.Range("I132").Value = 8420.833000000001
$ws.Range("J132").Value = 13472.8
$ws.Range("K132").Value = 25262.499
$ws.Range("L132").Value = 40418.39999999999
$ws.Range("M132").Value = -22732.499
$ws.Range("N132").Value = -45478.39999999999

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 6196.8887
$ws.Range("J74").Value = 3755.6
$ws.Range("L74").Value = 3755.6
$ws.Range("N74").Value = -5627.6
$ws.Range("H75").Value = 169982.86
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").Value = ""
$ws.Range("H77").Value = 6196.8887
$ws.Range("J77").Value = 3755.6
$ws.Range("L77").Value = 11266.8
$ws.Range("N77").Value = -20626.8
$ws.Range("H78").Value = 169982.86
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").Value = ""
$ws.Range("H122").Value = 2742.7273
$ws.Range("I122").Value = 2536.074
$ws.Range("K122").Value = 7608.222
$ws.Range("M122").Value = -5158.222
$ws.Range("H132").Value = 5899.8
$ws.Range("I132").Value = 2874.75
$ws.Range("K132").Value = 8624.25
$ws.Range("M132").Value = -6094.25
$ws.Range("H136").Value = 3063.6726
$ws.Range("I136").Value = 2643.25
$ws.Range("J136").Value = 4745.364
$ws.Range("K136").Value = 7929.75
$ws.Range("L136").Value = 14236.092
$ws.Range("M136").Value = -5379.75
$ws.Range("N136").Value = -19336.092

